$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New skill rows data: A (id), B (skill name), C (type), D (mana), E (cooldown)
$rows = @(
    @(36, "BabyDragonRoar", "PassiveSkill(Agumon)", 0, 0),
    @(37, "PepperBreath", "DamageSkill(Agumon)", 5, 2),
    @(38, "KnucklePunch", "DamageSkill(Monodramon)", 5, 2),
    @(39, "Dragon'sGrumble", "PassiveSkill(Veemon)", 0, 0),
    @(40, "VeemonHeadButt", "DamageSkill(Veemon)", 5, 2),
    @(41, "RelentlessAgression", "PassiveSkill(Guilmon)", 0, 0),
    @(42, "FireBall", "DamageSkill(Guilmon)", 5, 2),
    @(43, "Keen", "PassiveSkill(Dorumon)", 0, 0),
    @(44, "MetalKannon", "DamageSkill(Dorumon)", 5, 2),
    @(45, "MysticScale", "PassiveSkill(Betamon)", 0, 0),
    @(46, "ElectricShock", "DamageSkill(Betamon)", 5, 2)
)

$startRow = 38
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $rowRange = $ws.Range("A$r" + ":E$r")
    $rowRange.Style = "Good"
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("E$r").Value = $data[4]
}
